$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("paper")

# Column B (Electricity) and Column C (Heat) updates per row.
# $null means "leave as empty/inline string cell" (no numeric value).
$updates = @(
    @{ Row = 2;  B = 3.6770716;  C = 10.4493092 }
    @{ Row = 3;  B = 3.8431654;  C = 21.0865803 }
    @{ Row = 4;  B = 7.8961658;  C = 24.9691487 }
    @{ Row = 5;  B = 2.7493126;  C = 3.9613209;  D = $null }
    @{ Row = 6;  B = 3.4504604;  C = 6.6732969 }
    @{ Row = 7;  C = 4.3620868 }
    @{ Row = 8;  C = 3.6701993 }
    @{ Row = 9;  C = 5.4564957 }
    @{ Row = 10; C = 6.5406227 }
    @{ Row = 11; B = 2.1554668;  C = 1.9343046 }
    @{ Row = 12; B = 3.2669837;  C = 5.8615834 }
    @{ Row = 13; B = -1.5730766; C = -0.9853803 }
    @{ Row = 14; B = 18.3928806; C = 16.0873723 }
    @{ Row = 15; B = 3.8713396;  C = 8.7933913 }
    @{ Row = 16; B = 3.3037158;  C = 5.7348456 }
    @{ Row = 17; B = 3.2516518;  C = 10.2305804; D = $null }
    @{ Row = 18; B = 2.8585616;  C = 6.4040938 }
    @{ Row = 19; B = 5.1163425;  C = 17.7594201 }
    @{ Row = 20; B = 5.5030878;  C = 8.700299899999999;  D = $null }
    @{ Row = 21; B = 3.0226693;  C = 7.4471963 }
    @{ Row = 22; B = 4.9562555;  C = 20.1476858 }
    @{ Row = 23; B = 6.2428234;  C = 17.6932619 }
    @{ Row = 24; C = 14.5233804 }
    @{ Row = 25; B = 14.0802571; C = 11.7756491; D = $null }
    @{ Row = 26; B = 11.4523993; C = 6.4765879 }
    @{ Row = 29; B = 0.6250278;  C = 2.4420417 }
    @{ Row = 31; B = 3.8198693;  C = 6.9261068 }
    @{ Row = 34; B = 2.6851441;  C = -6.8902 }
    @{ Row = 35; B = 21.0227329; C = 13.0235623 }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey('B')) {
        $ws.Cells.Item($r, 2).Value = $u.B
    }
    if ($u.ContainsKey('C')) {
        $ws.Cells.Item($r, 3).Value = $u.C
    }
    if ($u.ContainsKey('D')) {
        $ws.Cells.Item($r, 4).Value = $u.D
    }
}
